# Insert a new weekly price record as row 85 (pushing the existing rows
# 85-156 down to 86-157), matching the "Fruta / hortaliza, semanal" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(85).Insert()

$ws.Range("A85").Value = 10
$ws.Range("B85").Value = "Vega Modelo de Temuco"
$ws.Range("C85").Value = "La Araucanía"
$ws.Range("D85").Value = 45040
$ws.Range("E85").Value = 9
$ws.Range("F85").Value = 100112035
$ws.Range("G85").Value = "Bruselas (repollito)"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 100
$ws.Range("K85").Value = 28000
$ws.Range("L85").Value = 28000
$ws.Range("M85").Value = 28000
$ws.Range("N85").Value = "$/malla 15 kilos"
$ws.Range("O85").Value = "Región Metropolitana"
$ws.Range("P85").Value = 1867
$ws.Range("Q85").Value = 15
$ws.Range("R85").Value = "Hortaliza"
